$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "Issue date" string (row 5, column A) ---
$ws.Range("A5").Value = "Issue date: 07/05/2021 17:42:28"

# --- Update column C width ---
$ws.Range("C1").ColumnWidth = 21.7109375

# --- Rewrite the results table (header row 15, data rows 16-21) ---

# Header row stays the same text, just re-written for safety
$ws.Range("A15").Value = "Number"
$ws.Range("B15").Value = "Language"
$ws.Range("C15").Value = "Features types"
$ws.Range("D15").Value = "Features selectors"
$ws.Range("E15").Value = "Pre-processing"
$ws.Range("F15").Value = "Technique"
$ws.Range("G15").Value = "MLP"
$ws.Range("H15").Value = "SVC"
$ws.Range("I15").Value = "LR"
$ws.Range("J15").Value = "RF"
$ws.Range("K15").Value = "MNB"
$ws.Range("L15").Value = "RNN"

# Clear old data rows 16-18 (old J column usage, etc.) and whole former table body
$ws.Range("A16:L18").Clear()

# Row 16
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "English"
$ws.Range("C16").Value = "TF: 10 words unigrams"
$ws.Range("D16").Value = "None"
$ws.Range("E16").Value = "None"
$ws.Range("F16").Value = "5 folds X 1 iterations CV"
$ws.Range("G16").Value = "66.5*"

# Row 17
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "English"
$ws.Range("C17").Value = "TF: 10 words unigrams"
$ws.Range("D17").Value = "None"
$ws.Range("E17").Value = "None"
$ws.Range("F17").Value = "5 folds X 1 iterations CV"
$ws.Range("G17").Value = "68.0*"

# Row 18
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = "english"
$ws.Range("C18").Value = "Doc2VecTransfomer"
$ws.Range("D18").Value = "None"
$ws.Range("E18").Value = "None"
$ws.Range("F18").Value = "5 folds X 20 iterations CV"
$ws.Range("L18").Value = "50.0*"

# Row 19
$ws.Range("A19").Value = 0
$ws.Range("B19").Value = "english"
$ws.Range("C19").Value = "Doc2VecTransfomer"
$ws.Range("D19").Value = "None"
$ws.Range("E19").Value = "None"
$ws.Range("F19").Value = "5 folds X 20 iterations CV"
$ws.Range("L19").Value = "76.98"

# Row 20
$ws.Range("A20").Value = 0
$ws.Range("B20").Value = "hebrew"
$ws.Range("C20").Value = "Doc2VecTransfomer"
$ws.Range("D20").Value = "None"
$ws.Range("E20").Value = "None"
$ws.Range("F20").Value = "5 folds X 20 iterations CV"
$ws.Range("L20").Value = "50.0*"

# Row 21
$ws.Range("A21").Value = 0
$ws.Range("B21").Value = "hebrew"
$ws.Range("C21").Value = "Doc2VecTransfomer"
$ws.Range("D21").Value = "None"
$ws.Range("E21").Value = "None"
$ws.Range("F21").Value = "5 folds X 20 iterations CV"
$ws.Range("L21").Value = "67.48*"

# Apply the "significantly smaller" style (red, style s6 in original) to the final RNN cell of row 19 (76.98)
# and keep others with the "centered wrap" plain cell style used across the table body.
# Re-apply font/number formatting consistent with original body style (font 8, centered, wrap text).
$bodyRange = $ws.Range("A16:L21")
$bodyRange.Font.Name = "Times New Roman"
$bodyRange.Font.Size = 10
$bodyRange.HorizontalAlignment = -4108  # xlCenter
$bodyRange.VerticalAlignment = -4108    # xlCenter
$bodyRange.WrapText = $true

# G17 (68.0*) uses the "significantly smaller" blue style (style 5: blue font)
$ws.Range("G17").Font.Color = 16711680   # Blue (BGR encoding for RGB 0000FF)
$ws.Range("G17").Font.Size = 10
$ws.Range("G17").HorizontalAlignment = -4108
$ws.Range("G17").VerticalAlignment = -4108

# L19 (76.98) uses the "significantly larger" red style (style 6: red font)
$ws.Range("L19").Font.Color = 255   # Red
$ws.Range("L19").Font.Size = 10
$ws.Range("L19").HorizontalAlignment = -4108
$ws.Range("L19").VerticalAlignment = -4108

# --- Resize / restyle the table (ListObject) to cover the new range ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A15:L21"))
$tbl.TableStyle = "TableStyleLight12"
